$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("C2").Value = 64
$ws.Range("D2").Value = 60.8
$ws.Range("E2").Value = 62.4
$ws.Range("F2").Value = "statistics"

# Update row 3 values
$ws.Range("B3").Value = "sairj.pdf"
$ws.Range("C3").Value = 22
$ws.Range("D3").Value = 24.4
$ws.Range("E3").Value = 23.2
$ws.Range("F3").Value = "python, statistics, pandas, data analysis"

# Remove rows 4 through 6 (no longer present in the updated sheet)
$ws.Range("A4:F6").Delete()
